$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "poke-ball"
$ws.Range("A3").Value = "poke-doll"
$ws.Range("A4").Value = "poke-radar"
$ws.Range("A5").Value = "slowpoke-tail"
$ws.Range("A6").Value = "pokeblock-case"
$ws.Range("A7").Value = "poke-flute"
$ws.Range("A8").Value = "poke-toy"
$ws.Range("A9").Value = "pokeblock-kit"
$ws.Range("A10").Value = "left-poke-ball"
$ws.Range("A11").Value = "pokemon-box"
